# Design-doc update:
#  - The cached "datetimeFigureOut" field text (shown on every slide layout,
#    the slide master, and the notes master via the Date placeholder) moves
#    from 1/2/2023 to 7/2/2023.
#  - Slide 5's "Customers whose ..." bullet drops the comma before "or need
#    flexible tool to integrate different apps."

function Set-FullRangeText {
    param($textRange, [string]$newText)
    # Replace the *entire* span as a single atomic edit (rather than letting
    # a plain `.Text =` assignment diff old/new and splinter the run at the
    # changed character), so we end up with one run carrying the full string.
    $len = $textRange.Text.Length
    if ($len -gt 0) {
        $full = $textRange.Characters(1, $len)
        $full.Text = $newText
    } else {
        $textRange.Text = $newText
    }
}

function Update-DatePlaceholder {
    param($shapes, [string]$oldDate, [string]$newDate)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                Set-FullRangeText $tr $newDate
            }
        }
    }
}

$p = $ppt.ActivePresentation
$oldDate = "1/2/2023"
$newDate = "7/2/2023"

# Slide master's own Date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes $oldDate $newDate

# Every slide layout's Date placeholder (Title Slide, Title and Content, ...).
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes $oldDate $newDate
}

# Notes master's Date placeholder.
if ($p.HasNotesMaster) {
    $notesMaster = $p.NotesMaster
    Update-DatePlaceholder $notesMaster.Shapes $oldDate $newDate
}

# Slide 5: "Customers whose system config setting and requirements keep
# updating, or need flexible tool..." -> drop the comma.
$slide5 = $p.Slides.Item(5)
for ($si = 1; $si -le $slide5.Shapes.Count; $si++) {
    $shape = $slide5.Shapes.Item($si)
    if ($shape.HasTextFrame) {
        $tf = $shape.TextFrame
        $tr = $tf.TextRange
        for ($pi = 1; $pi -le $tr.Paragraphs().Count; $pi++) {
            $para = $tr.Paragraphs($pi)
            if ($para.Text -eq "Customers whose system config setting and requirements keep updating, or need flexible tool to integrate different apps. ") {
                Set-FullRangeText $para "Customers whose system config setting and requirements keep updating or need flexible tool to integrate different apps. "
            }
        }
    }
}
